# Scenarios.xlsx: add "SteadyStateTime" / "SteadyStateTimeUnit" columns so the
# steady-state simulation time can be read from the scenario table.
# (commit: "Steady state time is read from Scenarios.xlsx  Fixes #323")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# The old "ModelFile" column lived at H. Insert two new blank columns at H:I
# so it slides over to J, then fill the two new header cells.
$ws.Columns("H:I").Insert()

$ws.Range("H1").Value = "SteadyStateTime"
$ws.Range("I1").Value = "SteadyStateTimeUnit"

# Row 2 (TestScenario / Indiv1) has no steady-state time - leave H2/I2 blank.

# Row 3 (TestScenario2 / Indiv) gets a steady state time of 500 minutes.
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = "min"

# Leave the selection where the author left it after editing the new columns.
$ws.Range("I4").Select() | Out-Null
